$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update header row (A1:D1) with new short English column codes
$ws.Range("A1").Value2 = "mx_state"
$ws.Range("B1").Value2 = "mx_municipality"
$ws.Range("C1").Value2 = "n_matriculas"
$ws.Range("D1").Value2 = "pct_matriculas"

# 2. Update "TOTAL" -> "Total" and capitalize lowercase connector words
#    (de/del/el/la/los/las/y) inside state / municipality name cells.
$ws.Range("B4").Value2 = 'Total'
$ws.Range("B8").Value2 = 'Total'
$ws.Range("B10").Value2 = 'Total'
$ws.Range("B13").Value2 = 'Total'
$ws.Range("B16").Value2 = 'Amatenango De La Frontera'
$ws.Range("B19").Value2 = 'Bejucal De Ocampo'
$ws.Range("B40").Value2 = 'Marqués De Comillas'
$ws.Range("B41").Value2 = 'Mazapa De Madero'
$ws.Range("B49").Value2 = 'San Cristóbal De Las Casas'
$ws.Range("B69").Value2 = 'Total'
$ws.Range("B75").Value2 = 'Total'
$ws.Range("A76").Value2 = 'Ciudad De México'
$ws.Range("B77").Value2 = 'Cuajimalpa De Morelos'
$ws.Range("B90").Value2 = 'Total'
$ws.Range("A91").Value2 = 'Coahuila De Zaragoza'
$ws.Range("B94").Value2 = 'Total'
$ws.Range("B104").Value2 = 'Total'
$ws.Range("A105").Value2 = 'Estado De México'
$ws.Range("B105").Value2 = 'Almoloya De Juárez'
$ws.Range("B106").Value2 = 'Atizapán De Zaragoza'
$ws.Range("B113").Value2 = 'Ecatepec De Morelos'
$ws.Range("B126").Value2 = 'Naucalpan De Juárez'
$ws.Range("B129").Value2 = 'San Felipe Del Progreso'
$ws.Range("B134").Value2 = 'Tenango Del Valle'
$ws.Range("B137").Value2 = 'Tlalnepantla De Baz'
$ws.Range("B143").Value2 = 'Total'
$ws.Range("B145").Value2 = 'Apaseo El Alto'
$ws.Range("B149").Value2 = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Range("B160").Value2 = 'San Francisco Del Rincón'
$ws.Range("B161").Value2 = 'San Luis De La Paz'
$ws.Range("B162").Value2 = 'San Miguel De Allende'
$ws.Range("B164").Value2 = 'Valle De Santiago'
$ws.Range("B167").Value2 = 'Total'
$ws.Range("B168").Value2 = 'Acapulco De Juárez'
$ws.Range("B171").Value2 = 'Ajuchitlán Del Progreso'
$ws.Range("B173").Value2 = 'Atlamajalcingo Del Monte'
$ws.Range("B175").Value2 = 'Atoyac De Álvarez'
$ws.Range("B176").Value2 = 'Ayutla De Los Libres'
$ws.Range("B179").Value2 = 'Chilapa De Álvarez'
$ws.Range("B180").Value2 = 'Chilpancingo De Los Bravo'
$ws.Range("B181").Value2 = 'Coahuayutla De José María Izazaga'
$ws.Range("B184").Value2 = 'Coyuca De Benítez'
$ws.Range("B185").Value2 = 'Coyuca De Catalán'
$ws.Range("B188").Value2 = 'Cuetzala Del Progreso'
$ws.Range("B189").Value2 = 'Cutzamala De Pinzón'
$ws.Range("B192").Value2 = 'Huitzuco De Los Figueroa'
$ws.Range("B193").Value2 = 'Iguala De La Independencia'
$ws.Range("B205").Value2 = 'Taxco De Alarcón'
$ws.Range("B206").Value2 = 'Tepecoacuilco De Trujano'
$ws.Range("B209").Value2 = 'Tlalixtaquilla De Maldonado'
$ws.Range("B210").Value2 = 'Tlapa De Comonfort'
$ws.Range("B215").Value2 = 'Zihuatanejo De Azueta'
$ws.Range("B216").Value2 = 'Total'
$ws.Range("B220").Value2 = 'Atotonilco De Tula'
$ws.Range("B221").Value2 = 'Atotonilco El Grande'
$ws.Range("B224").Value2 = 'Cuautepec De Hinojosa'
$ws.Range("B226").Value2 = 'Huejutla De Reyes'
$ws.Range("B229").Value2 = 'Jacala De Ledezma'
$ws.Range("B233").Value2 = 'Molango De Escamilla'
$ws.Range("B234").Value2 = 'Pachuca De Soto'
$ws.Range("B236").Value2 = 'Progreso De Obregón'
$ws.Range("B240").Value2 = 'Tepehuacán De Guerrero'
$ws.Range("B243").Value2 = 'Tulancingo De Bravo'
$ws.Range("B245").Value2 = 'Zacualtipán De Ángeles'
$ws.Range("B247").Value2 = 'Total'
$ws.Range("B248").Value2 = 'Ahualulco De Mercado'
$ws.Range("B249").Value2 = 'Encarnación De Díaz'
$ws.Range("B252").Value2 = 'Lagos De Moreno'
$ws.Range("B253").Value2 = 'Ojuelos De Jalisco'
$ws.Range("B255").Value2 = 'San Cristóbal De La Barranca'
$ws.Range("B257").Value2 = 'Tepatitlán De Morelos'
$ws.Range("B258").Value2 = 'Tizapán El Alto'
$ws.Range("B262").Value2 = 'Total'
$ws.Range("A263").Value2 = 'Michoacán De Ocampo'
$ws.Range("B286").Value2 = 'Total'
$ws.Range("B295").Value2 = 'Zacualpan De Amilpas'
$ws.Range("B296").Value2 = 'Total'
$ws.Range("B298").Value2 = 'Ixtlán Del Río'
$ws.Range("B300").Value2 = 'Total'
$ws.Range("B303").Value2 = 'Total'
$ws.Range("B307").Value2 = 'Cuyamecalco Villa De Zaragoza'
$ws.Range("B309").Value2 = 'Fresnillo De Trujano'
$ws.Range("B310").Value2 = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range("B311").Value2 = 'Heroica Ciudad De Tlaxiaco'
$ws.Range("B312").Value2 = 'Huajuapan De León'
$ws.Range("B313").Value2 = 'Ixtlán De Juárez'
$ws.Range("B315").Value2 = 'Mariscala De Juárez'
$ws.Range("B317").Value2 = 'Miahuatlán De Porfirio Díaz'
$ws.Range("B318").Value2 = 'Oaxaca De Juárez'
$ws.Range("B319").Value2 = 'Ocotlán De Morelos'
$ws.Range("B320").Value2 = 'Pinotepa De Don Luis'
$ws.Range("B321").Value2 = 'Putla Villa De Guerrero'
$ws.Range("B332").Value2 = 'San Francisco Del Mar'
$ws.Range("B334").Value2 = 'San José Del Progreso'
$ws.Range("B380").Value2 = 'Tezoatlán De Segura Y Luna'
$ws.Range("B381").Value2 = 'Villa De Tututepec'
$ws.Range("B382").Value2 = 'Villa De Tututepec De Melchor Ocampo'
$ws.Range("B383").Value2 = 'Villa De Zaachila'
$ws.Range("B384").Value2 = 'Total'
$ws.Range("B390").Value2 = 'Ayotoxco De Guerrero'
$ws.Range("B392").Value2 = 'Chalchicomula De Sesma'
$ws.Range("B398").Value2 = 'Chila De La Sal'
$ws.Range("B403").Value2 = 'Cuetzalan Del Progreso'
$ws.Range("B415").Value2 = 'Izúcar De Matamoros'
$ws.Range("B431").Value2 = 'San Salvador El Verde'
$ws.Range("B439").Value2 = 'Tepanco De López'
$ws.Range("B441").Value2 = 'Tepexi De Rodríguez'
$ws.Range("B442").Value2 = 'Tetela De Ocampo'
$ws.Range("B445").Value2 = 'Tlacotepec De Benito Juárez'
$ws.Range("B462").Value2 = 'Total'
$ws.Range("B466").Value2 = 'San Juan Del Río'
$ws.Range("B468").Value2 = 'Total'
$ws.Range("B471").Value2 = 'Total'
$ws.Range("B472").Value2 = 'Axtla De Terrazas'
$ws.Range("B476").Value2 = 'Ciudad Del Maíz'
$ws.Range("B481").Value2 = 'Santa María Del Río'
$ws.Range("B482").Value2 = 'Soledad De Graciano Sánchez'
$ws.Range("B484").Value2 = 'Villa De Guadalupe'
$ws.Range("B486").Value2 = 'Total'
$ws.Range("B491").Value2 = 'Total'
$ws.Range("B496").Value2 = 'Jalpa De Méndez'
$ws.Range("B501").Value2 = 'Total'
$ws.Range("B512").Value2 = 'Total'
$ws.Range("B516").Value2 = 'Ixtacuixtla De Mariano Matamoros'
$ws.Range("B518").Value2 = 'Papalotla De Xicohténcatl'
$ws.Range("B525").Value2 = 'Total'
$ws.Range("A526").Value2 = 'Veracruz De Ignacio De La Llave'
$ws.Range("B529").Value2 = 'Amatlán De Los Reyes'
$ws.Range("B537").Value2 = 'Cosamaloapan De Carpio'
$ws.Range("B538").Value2 = 'Cosautlán De Carvajal'
$ws.Range("B544").Value2 = 'Hueyapan De Ocampo'
$ws.Range("B545").Value2 = 'Ixhuacán De Los Reyes'
$ws.Range("B546").Value2 = 'Ixhuatlán De Madero'
$ws.Range("B547").Value2 = 'Ixhuatlán Del Café'
$ws.Range("B553").Value2 = 'Lerdo De Tejada'
$ws.Range("B557").Value2 = 'Paso De Ovejas'
$ws.Range("B579").Value2 = 'Total'
$ws.Range("B582").Value2 = 'Total'
$ws.Range("B585").Value2 = 'Nochistlán De Mejía'
$ws.Range("B590").Value2 = 'Total'
$ws.Range("A591").Value2 = 'Total'

# 3. Remove the trailing footnote / source rows (593-597); row 592 was
#    already an empty gap row so the data block ends at row 591.
$ws.Range("A593:A597").EntireRow.Delete()
